# Pflichtenheft aktualisiert - Projektkosten aktualisiert
#
# Updates the "interne Leistungen" (internal services) task descriptions
# in column C to reflect the more detailed / expanded scope of work, and
# adjusts row formatting (height + vertical alignment) so the longer
# descriptions remain readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content updates: expanded task descriptions (column C, "Beschreibung") ---

# Web-Entwickler (row 14)
$ws.Range("C14").Value = "25h, Logik, Ausgabe, Eingabe, Datenbankzugriff, Login(Verschlüsselung), "

# Datenbank-Entwickler (row 15)
$ws.Range("C15").Value = "8h, DB-Design, Erstellung"

# Projektleiter (row 16)
$ws.Range("C16").Value = "40h, Konzeption, Dokumenation, Aufgabenverteilung, Budgetverwaltung, Meilensteine festlegen"

# Berater (row 18) - Web-Designer (row 17) description stays unchanged
$ws.Range("C18").Value = "32h, Machbarkeit, Umfeld, Risiko, Statusberichte, Kommunikationsrichtlinien, Dokumentationsrichtlinien"

# --- Formatting updates ---

# Column C on rows 16/18 now wraps like row 14 already did
$ws.Range("C16").WrapText = $true
$ws.Range("C18").WrapText = $true

# Rows with the longer descriptions get taller to fit the wrapped text
$ws.Rows("14:14").RowHeight = 30
$ws.Rows("16:16").RowHeight = 30
$ws.Rows("18:18").RowHeight = 30

# Keep label/price cells aligned to the top of the now-taller rows
$ws.Range("A14:B14").VerticalAlignment = -4160
$ws.Range("A16:B16").VerticalAlignment = -4160
$ws.Range("A18:B18").VerticalAlignment = -4160

# --- View state (cursor / scroll position when the file was last saved) ---
$ws.Activate()
$ws.Range("F15").Select()
